$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = '???'
$ws.Range("F3").Value = '${begin.german}'
$ws.Range("G3").Value = '${ende.german}'
$ws.Range("D2").Value = '${bezeichnung}'

$ws.Range("D2:G2").Select
